# "Se sube caso de prueba" - remove the "Alta Deportista - Cancelar carga"
# test case row from the Test_Case sheet (entire row 6), shifting the rows
# below it up, just like deleting a row in the Excel UI.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_Case")
$ws.Rows.Item(6).Delete()
